$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the server data row (row 2).
# Shared-string insertion order matters for matching the target file's
# sharedStrings.xml index ordering, so we write A2, then C2, then B2,
# before the numeric cells.
$ws.Range("A2").Value = "SqlServer_1"
$ws.Range("C2").Value = "127.0.0.1"
$ws.Range("B2").Value = "000107001"
$ws.Range("D2").Value = 7001
$ws.Range("E2").Value = 123456

# Update the active selection to E4 (as left by the author after editing).
$ws.Range("E4").Select()
